# Update Lgals1-Cd69 LR-pair sheet with newly recomputed TPM values.
# The refreshed data only contains rows where the Target cluster is MuSCs,
# so the table shrinks from 6 data rows (2-7) down to 3 data rows (2-4),
# and the per-row numeric columns are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three now-obsolete data rows (old rows 5, 6 and 7) so the
# remaining three data rows shift up into rows 2-4.
$ws.Rows("5:7").Delete()

# Row 2: ECs -> Lgals1 -> Cd69 -> MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgals1"
$ws.Range("C2").Value = "Cd69"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.971347333333334
$ws.Range("H2").Value = 8.914042
$ws.Range("I2").Value = 0.02922956310646057
$ws.Range("J2").Value = 0.02922956310646057
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.103879
$ws.Range("N2").Value = 0.311637
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3086605896393333
$ws.Range("R2").Value = 2.777945306754
$ws.Range("S2").Value = 0.02922956310646057
$ws.Range("T2").Value = 0.02922956310646057

# Row 3: FAPs -> Lgals1 -> Cd69 -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgals1"
$ws.Range("C3").Value = "Cd69"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 63.88336466666667
$ws.Range("H3").Value = 191.650094
$ws.Range("I3").Value = 0.6284296749927923
$ws.Range("J3").Value = 0.6284296749927923
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.103879
$ws.Range("N3").Value = 0.311637
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 6.636140038208667
$ws.Range("R3").Value = 59.72526034387801
$ws.Range("S3").Value = 0.6284296749927923
$ws.Range("T3").Value = 0.6284296749927923

# Row 4: MuSCs -> Lgals1 -> Cd69 -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Lgals1"
$ws.Range("C4").Value = "Cd69"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 34.80083866666666
$ws.Range("H4").Value = 104.402516
$ws.Range("I4").Value = 0.3423407619007471
$ws.Range("J4").Value = 0.3423407619007471
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.103879
$ws.Range("N4").Value = 0.311637
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 3.615076319854666
$ws.Range("R4").Value = 32.535686878692
$ws.Range("S4").Value = 0.3423407619007471
$ws.Range("T4").Value = 0.3423407619007471
